$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "6927    EDMONTON, AB         T5Z3S2    164 AVE  "
$ws.Range("A5").Select()
